# Fix rPr child-element ordering in the custom "Tok" character styles so
# that <w:b/>/<w:i/> precede <w:color/>, matching the CT_RPr sequence in
# wml.xsd. OOXMLValidatorCLI flagged the previous order
# (<w:color/> then <w:b/>) as invalid ("unexpected child element"), even
# though xmllint didn't complain.
#
# Re-assigning Font.Bold / Font.Italic to their own current value is a
# functional no-op when that property is already switched on, but it
# marks the style dirty so its <w:rPr> gets rewritten in schema-canonical
# element order on save. (Only touch the property when it is already
# True, so we don't materialize a new explicit "off" element for styles
# that never carried <w:b/> / <w:i/> in the first place.)

$d = $word.ActiveDocument

$styleNames = @(
    "KeywordTok",
    "ImportTok",
    "CommentTok",
    "DocumentationTok",
    "AnnotationTok",
    "CommentVarTok",
    "ControlFlowTok",
    "InformationTok",
    "WarningTok",
    "AlertTok",
    "ErrorTok"
)

foreach ($name in $styleNames) {
    $s = $d.Styles($name)
    if ($s.Font.Bold) {
        $s.Font.Bold = $true
    }
    if ($s.Font.Italic) {
        $s.Font.Italic = $true
    }
}
